# NewlineInFormulas.xlsx — add a Fibonacci-style column on Sheet1:
#   B1 = 1, B2 = 2
#   B3:B10 = previous two cells summed (B3=B1+B2, B4=B2+B3, ... B10=B8+B9)
# and leave the selection on B3:B10 (active cell B3), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed values
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2

# Fill the whole range at once with a relative (R1C1) formula so every
# row gets its own correctly-offset relative formula:
#   B3  -> =B1+B2
#   B4  -> =B2+B3
#   ...
#   B10 -> =B8+B9
$ws.Range("B3:B10").FormulaR1C1 = "=R[-2]C+R[-1]C"

# Match the resulting selection state recorded in the sheet view.
$ws.Range("B3:B10").Select()
